# Update NATMI Wnt5b-Fzd2 LR-pairs sheet with newly-recomputed TPM numbers.
#
# The source data dropped the "Neutrophils" and "Resolving-Mac" target
# clusters entirely (for both the FAPs and MuSCs sending clusters), and the
# remaining rows were recomputed against the new TPM matrix. Net effect on
# the sheet: 13 data+header rows -> 9, and a handful of specificity /
# weight columns (K,L,M..T) get refreshed numbers in every surviving row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the rows that reference the dropped target clusters.
#    Delete from the bottom up so earlier row numbers stay valid.
#    Row 12-13: MuSCs -> Neutrophils / Resolving-Mac
#    Row 6-7  : FAPs  -> Neutrophils / Resolving-Mac
$ws.Rows("12:13").Delete()
$ws.Rows("6:7").Delete()

# After the deletions the sheet is laid out as:
#   row 2: FAPs  -> ECs
#   row 3: FAPs  -> FAPs
#   row 4: FAPs  -> Inflammatory-Mac
#   row 5: FAPs  -> MuSCs
#   row 6: MuSCs -> ECs
#   row 7: MuSCs -> FAPs
#   row 8: MuSCs -> Inflammatory-Mac
#   row 9: MuSCs -> MuSCs

# 2) Refresh the recomputed numeric columns (K,L are Receptor-expressing
#    cells / detection-rate counts; M..T are the recomputed expression /
#    specificity / weight values) for every surviving row.

$ws.Range("M2").Value = 0.1557005
$ws.Range("N2").Value = 0.311401
$ws.Range("O2").Value = 0.02102398211576467
$ws.Range("P2").Value = 0.01500040222529337
$ws.Range("Q2").Value = 0.05812984747199999
$ws.Range("R2").Value = 0.3487790848319999
$ws.Range("S2").Value = 0.01541175810434764
$ws.Range("T2").Value = 0.01207015424888988

$ws.Range("O3").Value = 0.7998659708565604
$ws.Range("P3").Value = 0.8560446272575798
$ws.Range("S3").Value = 0.5863466202959204
$ws.Range("T3").Value = 0.6888209089160175

$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.024117
$ws.Range("N4").Value = 0.072351
$ws.Range("O4").Value = 0.003256478795417461
$ws.Range("P4").Value = 0.003485197868350457
$ws.Range("Q4").Value = 0.009003937247999999
$ws.Range("R4").Value = 0.08103543523199999
$ws.Range("S4").Value = 0.002387181609580907
$ws.Range("T4").Value = 0.0028043831910027

$ws.Range("M5").Value = 1.3023455
$ws.Range("N5").Value = 2.604691
$ws.Range("O5").Value = 0.1758535682322574
$ws.Range("P5").Value = 0.1254697726487764
$ws.Range("Q5").Value = 0.4862228783519999
$ws.Range("R5").Value = 2.917337270112
$ws.Range("S5").Value = 0.1289105289596737
$ws.Range("T5").Value = 0.1009599267205155

$ws.Range("M6").Value = 0.1557005
$ws.Range("N6").Value = 0.311401
$ws.Range("O6").Value = 0.02102398211576467
$ws.Range("P6").Value = 0.01500040222529337
$ws.Range("Q6").Value = 0.021168105777
$ws.Range("R6").Value = 0.08467242310799999
$ws.Range("S6").Value = 0.005612224011417028
$ws.Range("T6").Value = 0.002930247976403487

$ws.Range("O7").Value = 0.7998659708565604
$ws.Range("P7").Value = 0.8560446272575798
$ws.Range("S7").Value = 0.21351935056064
$ws.Range("T7").Value = 0.1672237183415621

$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.024117
$ws.Range("N8").Value = 0.072351
$ws.Range("O8").Value = 0.003256478795417461
$ws.Range("P8").Value = 0.003485197868350457
$ws.Range("Q8").Value = 0.003278802618
$ws.Range("R8").Value = 0.019672815708
$ws.Range("S8").Value = 0.0008692971858365546
$ws.Range("T8").Value = 0.0006808146773477564

$ws.Range("M9").Value = 1.3023455
$ws.Range("N9").Value = 2.604691
$ws.Range("O9").Value = 0.1758535682322574
$ws.Range("P9").Value = 0.1254697726487764
$ws.Range("Q9").Value = 0.177059080107
$ws.Range("R9").Value = 0.7082363204279999
$ws.Range("S9").Value = 0.04694303927258368
$ws.Range("T9").Value = 0.02450984592826091
